$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$ws.Range("A2:A9").Value = "2025-10-05 18:29:01"
